$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column S (year 2022) mirrors the formatting of column R.
$ws.Range("R4").Copy($ws.Range("S4"))
$ws.Range("R5").Copy($ws.Range("S5"))

# Set the actual values for the new/changed cells.
$ws.Range("S4").Value = 2022
$ws.Range("Q5").Value = 91.892815141492093
$ws.Range("R5").Value = 101.53074848578628
$ws.Range("S5").Value = 109.27053140096621

# Q5 picks up the same number format / borders as R5 (style changed from 10 to 22).
$ws.Range("R5").Copy($ws.Range("Q5"))
$ws.Range("Q5").Value = 91.892815141492093

$ws.Range("T5").Select() | Out-Null
